# This workbook holds a weekly rolling series of price observations (rows 2-18).
# On each weekly refresh the observations are re-shuffled among the existing
# rows (same 17 data rows, values rotated to different row positions).
# The mapping below says: new row <key> gets the data that currently (before
# this edit) lives in row <value>, for the columns that vary week to week
# (Fecha, Calidad, Volumen, Precio mínimo/máximo/promedio, Unidad de
# comercialización, Origen, Precio $/Kg, Kg / unidad).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 10
    3  = 5
    4  = 6
    5  = 11
    6  = 7
    7  = 13
    8  = 14
    9  = 3
    10 = 18
    11 = 15
    12 = 16
    13 = 17
    14 = 9
    15 = 4
    16 = 12
    17 = 2
    18 = 8
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot current (pre-edit) values for every relevant cell before writing
# anything, since several rows exchange values with each other.
# NOTE: use .Value2 (not .Value) when *reading* - in this runtime reading
# .Value back out into a variable does not reliably yield the underlying
# scalar, while .Value2 does (for both numbers and strings).
$snapshot = @{}
foreach ($row in 2..18) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowVals
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $sourceVals = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $sourceVals[$col]
    }
}
